# docs/questions/qs-pmfspdfscdfs.docx -- "new guides and refresh"
#
# 1. Collapse the word-by-word runs in the Title / Author / Abstract
#    paragraphs into single runs holding the full sentence.
# 2. Strip the redundant explicit "left" paragraph justification
#    (<w:jc w:val="left"/>) that was being written on every "Compact"
#    style paragraph (mostly inside the quiz tables) -- left is already
#    the inherited default, so clearing/re-asserting alignment removes
#    the now-superfluous attribute.

$d = $word.ActiveDocument

function Merge-ParagraphRuns($para, [string]$fullText) {
    $rng = $para.Range
    # Trim the trailing paragraph mark from the comparison/search range
    # so Find only matches the visible text of this paragraph.
    $rng.MoveEnd(1, -1) | Out-Null
    $rng.Find.ClearFormatting()
    $rng.Find.Execute($fullText, $true, $false, $false, $false, $false, $true, 1, $false, $fullText, 2) | Out-Null
}

Merge-ParagraphRuns $d.Paragraphs(1) "Questions: PMFs, PDFs, and CDFs"
Merge-ParagraphRuns $d.Paragraphs(2) "Sophie Chowgule"
Merge-ParagraphRuns $d.Paragraphs(4) "A selection of questions to test your understanding of probability mass functions (PMFs), probability density functions (PDFs), and cumulative distribution functions (CDFs)."

# Re-assert (wdAlignParagraphLeft = 0) on every paragraph in the body so
# any explicit <w:jc w:val="left"/> collapses back to the inherited
# default and is dropped from the saved markup.
foreach ($p in $d.Paragraphs) {
    $p.Alignment = 0
}

Write-Output "edit complete"
